$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.378.18'
$ws.Range("E2").Value = '  +1.48%  '
$ws.Range("D3").Value = '1.622.85'
$ws.Range("E3").Value = '  +1.94%  '
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.27'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.52%  '
$ws.Range("E6").Value = '  -0.20%  '
$ws.Range("E7").Value = '  +0.92%  '
$ws.Range("E8").Value = '  +0.50%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0616'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.61%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.88'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +4.53%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0815'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.64%  '
$ws.Range("D12").Value = '1.848.76'
$ws.Range("E12").Value = '  +1.99%  '
$ws.Range("D13").Value = '1.619.75'
$ws.Range("E13").Value = '  +1.68%  '
$ws.Range("E14").Value = '  +0.43%  '
$ws.Range("E15").Value = '  +1.18%  '
$ws.Range("D16").Value = '26.360.86'
$ws.Range("E16").Value = '  +1.48%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.57'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +3.14%  '
$ws.Range("E18").Value = '  +0.44%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '202.45'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +0.08%  '
$ws.Range("E21").Value = '  +0.74%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.35'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +1.35%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.05'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.79%  '
$ws.Range("E24").Value = '  -1.82%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.70'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.85%  '
$ws.Range("E26").Value = '  -0.15%  '
$ws.Range("E27").Value = '  -0.99%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.19'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.29%  '
$ws.Range("E29").Value = '  +1.36%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0519'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +9.44%  '
$ws.Range("E31").Value = '  +0.52%  '
$ws.Range("E32").Value = '  +2.03%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.92'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +0.93%  '
$ws.Range("E34").Value = '  +1.63%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.39'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +1.82%  '
$ws.Range("D36").Value = '1.175.47'
$ws.Range("E36").Value = '  +4.15%  '
$ws.Range("E37").Value = '  +0.53%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.808'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +1.85%  '
$ws.Range("E39").Value = '  -0.18%  '
$ws.Range("E40").Value = '  +0.10%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.497'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +1.65%  '
$ws.Range("E42").Value = '  +4.99%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.786'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +0.92%  '
$ws.Range("D44").Value = '1.760.62'
$ws.Range("E44").Value = '  +2.09%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '92.71'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.53%  '
$ws.Range("E46").Value = '  +2.99%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '53.97'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +0.53%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0509'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.77%  '
$ws.Range("E49").Value = '  +0.68%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.00'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -0.60%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.28'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +1.75%  '
